$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "Condicion_Pacientes" currently spans A1:F41 (header + 40 data rows).
# Add a new row to the table; this automatically grows the table ref,
# the autoFilter ref, and the worksheet dimension.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Copy the formatting (date number format) from the cell above so the new
# date cell keeps the same style index instead of creating a new one.
$ws.Range("A41").Copy($ws.Range("A42"))

# Fill in the new data row.
$ws.Range("A42").Value = 43961
$ws.Range("B42").Value = 432
$ws.Range("C42").Value = 142
$ws.Range("D42").Value = 236
$ws.Range("E42").Value = 6
$ws.Range("F42").Value = 11

# Move the active selection, matching the state left behind by the edit.
$ws.Range("F44").Select() | Out-Null
